# "added experimental to vs"
#
# 1. Fix the sheet-name typo: "Include from SNOWMED CT" -> "Include from SNOMED CT"
# 2. Populate the "Experimental" row (B7) on the Metadata sheet with the
#    literal text value "true" (the Status/Experimental flag for this ValueSet).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from SNOWMED CT")

# --- Fix sheet name typo ---------------------------------------------------
$ws2.Name = "Include from SNOMED CT"

# --- Set Experimental = true (as literal text, not boolean) ----------------
# Writing the bare word true/false to a cell auto-converts it to a Boolean,
# so build it as a text formula first and then convert the cell to a plain
# value in place (Copy + PasteSpecial values-only), which keeps the cell's
# existing style/format and stores it as literal text "true".
$cell = $ws1.Range("B7")
$cell.Formula = "=""true"""
$cell.Copy()
$cell.PasteSpecial(-4163)

# --- Refresh the generated "Date" metadata timestamp ------------------------
$ws1.Range("B8").Value = "2024-05-23T15:19:20+10:00"
